# Generate Report for Handoff
# Updates the "Latest Handback DateTime" column (D) for the most recently
# handed-back file (row 7, bb3e754f-...) on both the zh-cn and de-de
# localization-status sheets with fresh handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-10 03:13:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-10 03:13:36"
